$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.578.40'
$ws.Range('D3').Value = '2.458.09'
$ws.Range('E3').Value = '  -0.98%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').Value = '314.58'
$ws.Range('E5').Value = '  +0.86%  '
$ws.Range('D6').Value = '92.15'
$ws.Range('E6').Value = '  -0.39%  '
$ws.Range('E7').Value = '  +1.86%  '
$ws.Range('E8').Value = '  -0.33%  '
$ws.Range('E9').Value = '  +3.78%  '
$ws.Range('D10').Value = '32.41'
$ws.Range('E10').Value = '  +1.23%  '
$ws.Range('D11').Value = '0.0795'
$ws.Range('E11').Value = '  +2.90%  '
$ws.Range('E12').Value = '  +1.00%  '
$ws.Range('D13').Value = '2.835.99'
$ws.Range('E13').Value = '  -0.41%  '
$ws.Range('D14').Value = '6.83'
$ws.Range('E14').Value = '  +1.00%  '
$ws.Range('E15').Value = '  +4.28%  '
$ws.Range('D16').Value = '2.514.08'
$ws.Range('E16').Value = '  +0.30%  '
$ws.Range('D17').Value = '0.773'
$ws.Range('E17').Value = '  +1.76%  '
$ws.Range('D18').Value = '41.582.02'
$ws.Range('E18').Value = '  +0.33%  '
$ws.Range('E19').Value = '  +3.37%  '
$ws.Range('D20').Value = '0.0₃0935'
$ws.Range('E20').Value = '  +2.18%  '
$ws.Range('D21').Value = '70.73'
$ws.Range('E21').Value = '  +0.46%  '
$ws.Range('D22').Value = '11.34'
$ws.Range('E22').Value = '  +2.89%  '
$ws.Range('D23').Value = '237.95'
$ws.Range('E23').Value = '  +1.70%  '
$ws.Range('D24').Value = '2.70'
$ws.Range('E24').Value = '  +0.41%  '
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('E26').Value = '  +1.47%  '
$ws.Range('D27').Value = '24.28'
$ws.Range('E27').Value = '  +0.24%  '
$ws.Range('E28').Value = '  +1.08%  '
$ws.Range('E29').Value = '  +0.92%  '
$ws.Range('E30').Value = '  -3.30%  '
$ws.Range('D31').Value = '155.90'
$ws.Range('E31').Value = '  +1.49%  '
$ws.Range('E32').Value = '  +2.13%  '
$ws.Range('D33').Value = '2.58'
$ws.Range('E33').Value = '  +0.95%  '
$ws.Range('D34').Value = '0.0758'
$ws.Range('E34').Value = '  +0.49%  '
$ws.Range('D35').Value = '17.43'
$ws.Range('E35').Value = '  -3.05%  '
$ws.Range('D36').Value = '2.44'
$ws.Range('E36').Value = '  -1.91%  '
$ws.Range('E37').Value = '  -2.90%  '
$ws.Range('E38').Value = '  +2.02%  '
$ws.Range('E39').Value = '  +1.73%  '
$ws.Range('D40').Value = '1.79'
$ws.Range('E40').Value = '  -1.78%  '
$ws.Range('E41').Value = '  -3.82%  '
$ws.Range('E42').Value = '  -0.46%  '
$ws.Range('D43').Value = '1.969.97'
$ws.Range('E43').Value = '  +1.24%  '
$ws.Range('E44').Value = '  +0.89%  '
$ws.Range('D45').Value = '18.61'
$ws.Range('E45').Value = '  -9.85%  '
$ws.Range('E46').Value = '  -1.06%  '
$ws.Range('E47').Value = '  +2.71%  '
$ws.Range('D48').Value = '2.694.84'
$ws.Range('E48').Value = '  -0.71%  '
$ws.Range('D49').Value = '96.34'
$ws.Range('E49').Value = '  +1.04%  '
$ws.Range('D50').Value = '66.64'
$ws.Range('E50').Value = '  +0.62%  '
$ws.Range('E51').Value = '  -1.92%  '
